$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Tipo de formación", "tipo_de_formacion"),
    @("Limitaciones", "limitaciones"),
    @("Tipo de limitación", "tipo_de_limitacion"),
    @("Necesidades básicas insatisfechas", "necesidades_basicas_insatisfechas")
)

$startRow = 60
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
